# Apply the "ploidy level" values edit to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: row height change 12.5 -> 12.65 ---
$ws.Rows.Item(2).RowHeight = 12.65

# --- Fill in ploidy-level values for columns E,F,G,H across rows 2-169 ---
# The data falls into five contiguous bands where (E,F,G,H) is constant per band.
$ws.Range("E2:E25").Value = -1
$ws.Range("F2:F25").Value = 1
$ws.Range("G2:G25").Value = -1
$ws.Range("H2:H25").Value = -1

$ws.Range("E26").Value = 1
$ws.Range("F26").Value = -1
$ws.Range("G26").Value = -1
$ws.Range("H26").Value = -1

$ws.Range("E27:E158").Value = -1
$ws.Range("F27:F158").Value = 1
$ws.Range("G27:G158").Value = -1
$ws.Range("H27:H158").Value = -1

$ws.Range("E159:E162").Value = -1
$ws.Range("F159:F162").Value = -1
$ws.Range("G159:G162").Value = 1
$ws.Range("H159:H162").Value = -1

$ws.Range("E163:E169").Value = -1
$ws.Range("F163:F169").Value = -1
$ws.Range("G163:G169").Value = -1
$ws.Range("H163:H169").Value = 1

# --- View state: selection moves to H3, top-left scrolls to B3 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 2
$ws.Range("H3").Select() | Out-Null
